$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data rows: recording script text + calculated frequency.
# Fill the data rows (B2:C18) before the header cells so shared-string
# indices land in the same order the original authoring tool produced them.
$data = @(
    @("eh nice cole dower", 259836),
    @("a nice coal dower", 7747545),
    @("on ice coal dower", 2806985),
    @("an ice kohl dower", 807375),
    @("an ice cole dower", 808536),
    @("an ice coal dower", 826911),
    @("an aye scold hour", 996702),
    @("a nye scold hour", 7608364),
    @("a nigh scold our", 8011331),
    @("on ice cold hour", 2911102),
    @("an ice-cold hour", 866031),
    @("an ice-cold hour", 866031),
    @("an eye scold our", 1294559),
    @("a nye scold our", 8009974),
    @("a nice cold our", 8253272),
    @("an ice-cold our", 1267641),
    @("an ice cold our", 1332638)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

# Column headers added last
$ws.Range("C1").Value = "calculated_Freq"
$ws.Range("B1").Value = "recording script"

$ws.Columns.Item(2).ColumnWidth = 18.6

$ws.Range("A16").Select()
